$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (current "IOT folder" row) to make room for
# the new "MARIO" import row, shifting the existing rows down.
$ws.Rows.Item(2).Insert()

# Fill in the new row 2 with the MARIO github import path.
$ws.Range("A2").Value = "MARIO"
$ws.Range("B2").Value = "C:\Users\loren\Documents\GitHub\Waste-MARIO\DWMRIO\mario"

# Update the selection to match the target state.
$ws.Range("B2").Select()
